# Updated symbol list on Fri Feb 17 18:18:31 UTC 2023 with GitHub Actions
# Refreshes Price (D), Volume(1h) (E), and Hora (G) columns for each coin row (rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).Value = "'311.44"
$ws.Cells.Item(2, 5).Value = "'-2.96%"
$ws.Cells.Item(2, 7).Value = "'18"

$ws.Cells.Item(3, 4).Value = "'53.41"
$ws.Cells.Item(3, 5).Value = "'8.73%"
$ws.Cells.Item(3, 7).Value = "'18"

$ws.Cells.Item(4, 4).Value = "'5.092"
$ws.Cells.Item(4, 5).Value = "'-4.10%"
$ws.Cells.Item(4, 7).Value = "'18"

$ws.Cells.Item(5, 4).Value = "'0.07928"
$ws.Cells.Item(5, 5).Value = "'-1.69%"
$ws.Cells.Item(5, 7).Value = "'18"

$ws.Cells.Item(6, 4).Value = "'4.566"
$ws.Cells.Item(6, 5).Value = "'-1.08%"
$ws.Cells.Item(6, 7).Value = "'18"

$ws.Cells.Item(7, 4).Value = "'1.394"
$ws.Cells.Item(7, 5).Value = "'4.35%"
$ws.Cells.Item(7, 7).Value = "'18"

$ws.Cells.Item(8, 4).Value = "'1.670"
$ws.Cells.Item(8, 5).Value = "'1.82%"
$ws.Cells.Item(8, 7).Value = "'18"

$ws.Cells.Item(9, 5).Value = "'-1.55%"
$ws.Cells.Item(9, 7).Value = "'18"

$ws.Cells.Item(10, 4).Value = "'0.2012"
$ws.Cells.Item(10, 5).Value = "'2.08%"
$ws.Cells.Item(10, 7).Value = "'18"

$ws.Cells.Item(11, 4).Value = "'0.09517"
$ws.Cells.Item(11, 5).Value = "'-1.71%"
$ws.Cells.Item(11, 7).Value = "'18"

$ws.Cells.Item(12, 4).Value = "'0.04723"
$ws.Cells.Item(12, 5).Value = "'0.91%"
$ws.Cells.Item(12, 7).Value = "'18"

$ws.Cells.Item(13, 4).Value = "'0.1042"
$ws.Cells.Item(13, 5).Value = "'-0.62%"
$ws.Cells.Item(13, 7).Value = "'18"

$ws.Cells.Item(14, 4).Value = "'0.001282"
$ws.Cells.Item(14, 5).Value = "'-3.37%"
$ws.Cells.Item(14, 7).Value = "'18"

$ws.Cells.Item(15, 4).Value = "'0.005888"
$ws.Cells.Item(15, 5).Value = "'1.30%"
$ws.Cells.Item(15, 7).Value = "'18"

$ws.Cells.Item(16, 4).Value = "'3.337"
$ws.Cells.Item(16, 5).Value = "'-0.25%"
$ws.Cells.Item(16, 7).Value = "'18"

$ws.Cells.Item(17, 4).Value = "'2.435"
$ws.Cells.Item(17, 5).Value = "'-0.35%"
$ws.Cells.Item(17, 7).Value = "'18"

$ws.Cells.Item(18, 4).Value = "'0.3477"
$ws.Cells.Item(18, 5).Value = "'-1.35%"
$ws.Cells.Item(18, 7).Value = "'18"

$ws.Cells.Item(19, 4).Value = "'8.417"
$ws.Cells.Item(19, 5).Value = "'5.10%"
$ws.Cells.Item(19, 7).Value = "'18"

$ws.Cells.Item(20, 5).Value = "'-0.15%"
$ws.Cells.Item(20, 7).Value = "'18"

$ws.Cells.Item(21, 4).Value = "'0.2912"
$ws.Cells.Item(21, 5).Value = "'-5.88%"
$ws.Cells.Item(21, 7).Value = "'18"

$ws.Cells.Item(22, 4).Value = "'0.04174"
$ws.Cells.Item(22, 5).Value = "'-0.85%"
$ws.Cells.Item(22, 7).Value = "'18"

$ws.Cells.Item(23, 4).Value = "'0.001259"
$ws.Cells.Item(23, 5).Value = "'-3.90%"
$ws.Cells.Item(23, 7).Value = "'18"

$ws.Cells.Item(24, 4).Value = "'0.003976"
$ws.Cells.Item(24, 5).Value = "'-8.61%"
$ws.Cells.Item(24, 7).Value = "'18"

$ws.Cells.Item(25, 5).Value = "'-0.22%"
$ws.Cells.Item(25, 7).Value = "'18"

$ws.Cells.Item(26, 4).Value = "'0.0003531"
$ws.Cells.Item(26, 5).Value = "'-0.23%"
$ws.Cells.Item(26, 7).Value = "'18"

$ws.Cells.Item(27, 7).Value = "'18"

$ws.Cells.Item(28, 7).Value = "'18"

$ws.Cells.Item(29, 7).Value = "'18"

$ws.Cells.Item(30, 7).Value = "'18"

$ws.Cells.Item(31, 7).Value = "'18"

$ws.Cells.Item(32, 7).Value = "'18"

$ws.Cells.Item(33, 7).Value = "'18"

$ws.Cells.Item(34, 7).Value = "'18"

$ws.Cells.Item(35, 7).Value = "'18"

$ws.Cells.Item(36, 7).Value = "'18"

$ws.Cells.Item(37, 7).Value = "'18"

$ws.Cells.Item(38, 4).Value = "'0.02633"
$ws.Cells.Item(38, 5).Value = "'-3.11%"
$ws.Cells.Item(38, 7).Value = "'18"

$ws.Cells.Item(39, 4).Value = "'0.05951"
$ws.Cells.Item(39, 5).Value = "'-1.11%"
$ws.Cells.Item(39, 7).Value = "'18"

$ws.Cells.Item(40, 4).Value = "'0.01082"
$ws.Cells.Item(40, 5).Value = "'-0.36%"
$ws.Cells.Item(40, 7).Value = "'18"

$ws.Cells.Item(41, 4).Value = "'0.1704"
$ws.Cells.Item(41, 5).Value = "'15.90%"
$ws.Cells.Item(41, 7).Value = "'18"

$ws.Cells.Item(42, 4).Value = "'0.007945"
$ws.Cells.Item(42, 5).Value = "'-1.22%"
$ws.Cells.Item(42, 7).Value = "'18"

$ws.Cells.Item(43, 4).Value = "'0.008179"
$ws.Cells.Item(43, 5).Value = "'3.41%"
$ws.Cells.Item(43, 7).Value = "'18"

$ws.Cells.Item(44, 4).Value = "'0.008353"
$ws.Cells.Item(44, 5).Value = "'5.92%"
$ws.Cells.Item(44, 7).Value = "'18"

$ws.Cells.Item(45, 4).Value = "'0.3430"
$ws.Cells.Item(45, 5).Value = "'-1.70%"
$ws.Cells.Item(45, 7).Value = "'18"

$ws.Cells.Item(46, 4).Value = "'0.00007236"
$ws.Cells.Item(46, 5).Value = "'4.76%"
$ws.Cells.Item(46, 7).Value = "'18"

$ws.Cells.Item(47, 5).Value = "'-0.22%"
$ws.Cells.Item(47, 7).Value = "'18"

$ws.Cells.Item(48, 4).Value = "'0.06440"
$ws.Cells.Item(48, 5).Value = "'9.72%"
$ws.Cells.Item(48, 7).Value = "'18"

$ws.Cells.Item(49, 4).Value = "'0.002613"
$ws.Cells.Item(49, 5).Value = "'-34.66%"
$ws.Cells.Item(49, 7).Value = "'18"

$ws.Cells.Item(50, 4).Value = "'0.00002095"
$ws.Cells.Item(50, 5).Value = "'-0.22%"
$ws.Cells.Item(50, 7).Value = "'18"

$ws.Cells.Item(51, 4).Value = "'0.0001995"
$ws.Cells.Item(51, 5).Value = "'-0.22%"
$ws.Cells.Item(51, 7).Value = "'18"
